$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.852.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.236.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  -4.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.801.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "67.903.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.254.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "396.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.513"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000117"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("E26").Value = "  +2.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "161.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("E36").Value = "  -4.57%  "
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.809"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0683"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.605.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "335.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.12%  "
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "30.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.90%  "
